# Update for release: refresh the build timestamp embedded in the version
# string everywhere it appears in the workbook (new build ran later the
# same day: 17.29.55 EST -> 18.05.36 EST on February 03 2026).
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 18.05.36 EST)"

$wb = $excel.ActiveWorkbook

# --- "About" sheet ---------------------------------------------------------
$aboutWs = $wb.Worksheets.Item("About")

$aboutWs.Range("A2").Value = "Version: " + $newVersion

$newCitation = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Ensham Coal Mine, Australia, M0038, version ''' + $newVersion + '''. (See the CC license for attribution requirements if sharing or adapting the data set.)'
$aboutWs.Range("A6").Value = $newCitation

# --- "Boundaries and methane sources" sheet ---------------------------------
$dataWs = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S ("build_version") holds the same version string for every data
# row (rows 2 through 9).
for ($row = 2; $row -le 9; $row++) {
    $dataWs.Range("S" + $row).Value = $newVersion
}
